# NIT-9002819067.xlsx update
# Commit: "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The worker data table (rows 16-42, columns B:G) is rebuilt:
#   - The single "AROLDO SANJUAN LIÑAN" (CC 73138476, periodo 1701) row,
#     which used to sit at the top (row 16), now moves to the bottom (row 42).
#   - The 26 "LUIS GABRIEL RODRIGUEZ DEVOZ" (CC 1051443518) rows, which used
#     to run from row 17 (periodo 1809) ascending through row 42 (periodo 2010),
#     now run from row 16 (periodo 2010) descending through row 41 (periodo 1809).
#   - Column B (Tipo Doc) stays "CC" throughout; F/G (Valor Mora / Salario
#     Basico) travel together with their owning periodo row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New row 16..41: LUIS GABRIEL RODRIGUEZ DEVOZ, periods 2010 -> 1809 (descending)
$luisPeriods = @(
    "2010","2009","2008","2007","2006","2005","2004","2003","2002","2001",
    "1912","1911","1910","1909","1908","1907","1906","1905","1904","1903",
    "1902","1901","1812","1811","1810","1809"
)
$luisMora = @(
    33333,40000,40000,40000,40000,40000,40000,40000,40000,40000,
    40000,40000,40000,40000,40000,40000,40000,40000,40000,40000,
    40000,40000,40000,40000,40000,14667
)
$luisSalario = @(
    1000000,1000000,1000000,1000000,1000000,1000000,1000000,1000000,1000000,1000000,
    1000000,1000000,1000000,1000000,1000000,1000000,1000000,1000000,1000000,1000000,
    1000000,1000000,1000000,1000000,1000000,1000000
)

for ($i = 0; $i -lt $luisPeriods.Length; $i++) {
    $r = 16 + $i
    $ws.Range("B$r").Value = "CC"
    $ws.Range("C$r").Value = "1051443518"
    $ws.Range("D$r").Value = "LUIS GABRIEL RODRIGUEZ DEVOZ"
    $ws.Range("E$r").Value = $luisPeriods[$i]
    $ws.Range("F$r").Value = $luisMora[$i]
    $ws.Range("G$r").Value = $luisSalario[$i]
}

# New row 42: AROLDO SANJUAN LIÑAN, periodo 1701 (moved from old row 16)
$ws.Range("B42").Value = "CC"
$ws.Range("C42").Value = "73138476"
$ws.Range("D42").Value = "AROLDO SANJUAN LIÑAN"
$ws.Range("E42").Value = "1701"
$ws.Range("F42").Value = 27578
$ws.Range("G42").Value = 689454
